$wb = $excel.ActiveWorkbook

# --- Sheet "User role": update per-role Uuid checksums (column FS) ---
$wsRole = $wb.Worksheets.Item("User role")
$wsRole.Range("FS2").Value = "SVBQ4Q-M6LXCK-EILBRB-RV7RKHEQ"
$wsRole.Range("FS3").Value = "T2YFVQ-UALTV7-C6FPVV-5XUUSBAQ"
$wsRole.Range("FS4").Value = "WXUVOC-2FEOLI-MBR6FK-R2PO2KM4"
$wsRole.Range("FS5").Value = "TD5YWN-CTNFSD-OXU5UQ-P7T72BOM"
$wsRole.Range("FS6").Value = "UYDITU-7Q6BHJ-7UB4GK-WXEUSCIE"
$wsRole.Range("FS7").Value = "TSNCUO-G52ZP2-A7N5T2-M4J2KBAE"
$wsRole.Range("FS8").Value = "QALKMH-U23NN7-WTJA7L-73E2CHKM"
$wsRole.Range("FS9").Value = "XNR6EZ-KJDB7Y-LI57LH-4VF42I2Y"
$wsRole.Range("FS10").Value = "SZ3DX5-4AGEJE-ZGPSYH-SVNFCFM4"
$wsRole.Range("FS11").Value = "V27GIG-USYWMJ-G2JPTP-4SHJCBY4"
$wsRole.Range("FS12").Value = "WWCW3B-A6E6X4-3JZ4YF-P5AKKKVI"
$wsRole.Range("FS13").Value = "TKDQAO-JJGR7Q-XJAV7V-V2RKCORE"
$wsRole.Range("FS14").Value = "QMRO57-QLHTNI-RI6J3T-6L4Q2PAM"
$wsRole.Range("FS15").Value = "VN7AL2-XTVNCR-DNHKZD-SG2NKCAQ"
$wsRole.Range("FS16").Value = "T4N6N4-U63KOE-MID5YP-LNO6CNV4"
$wsRole.Range("FS17").Value = "SAPEKQ-BKXSWW-WGC4DO-N4YTSOXA"
$wsRole.Range("FS18").Value = "TA6WR3-PDRG4C-YQGWXI-I4CISIAI"
$wsRole.Range("FS19").Value = "XL3EZD-QZMGLW-H5RC3V-DN3BCGDE"
$wsRole.Range("FS20").Value = "WA7ZUG-HOILVH-TBPP7C-XSRQSCRQ"
$wsRole.Range("FS21").Value = "TRBFTM-MGJ2EP-OOSZG4-B3J7KJDE"
$wsRole.Range("FS22").Value = "UYIJ57-TIVXP3-GZUGCT-TIOM2PDU"
$wsRole.Range("FS23").Value = "RVZ4CV-MTFUKQ-C6A5HV-IJHZCIGM"
$wsRole.Range("FS24").Value = "VWIIWC-6YFPZI-HLP4PW-C7MMCEFY"
$wsRole.Range("FS25").Value = "RTO2PS-2MTJJH-RMBU4T-FUYUKBHA"
$wsRole.Range("FS26").Value = "RENBVH-NKED45-67G6VS-FB5QCDB4"
$wsRole.Range("FS27").Value = "UVM2VS-WYGJMH-VW3KPF-XJUE2FWQ"

# --- Sheet "User Rights": reorder "Needed user rights" lists (column E) ---
$wsRights = $wb.Worksheets.Item("User Rights")
$wsRights.Range("E4").Value = "PERSON_EDIT, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E6").Value = "DOCUMENT_VIEW, VISIT_DELETE, TASK_DELETE, IMMUNIZATION_DELETE, TREATMENT_DELETE, PATHOGEN_TEST_DELETE, CLINICAL_COURSE_VIEW, CLINICAL_VISIT_DELETE, PERSON_VIEW, SAMPLE_DELETE, ADDITIONAL_TEST_DELETE, PERSON_DELETE, TASK_VIEW, IMMUNIZATION_VIEW, PRESCRIPTION_DELETE, THERAPY_VIEW, DOCUMENT_DELETE, CASE_VIEW, SAMPLE_VIEW, ADDITIONAL_TEST_VIEW"
$wsRights.Range("E9").Value = "CASE_EDIT, PERSON_EDIT, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E10").Value = "CASE_EDIT, PERSON_EDIT, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E11").Value = "CASE_EDIT, PERSON_EDIT, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E12").Value = "CASE_EDIT, PERSON_EDIT, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E13").Value = "CASE_EDIT, PERSON_EDIT, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E14").Value = "CASE_EDIT, PERSON_EDIT, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E15").Value = "CASE_EDIT, PERSON_EDIT, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E17").Value = "CASE_EDIT, PERSON_EDIT, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E19").Value = "IMMUNIZATION_VIEW, PERSON_VIEW"
$wsRights.Range("E20").Value = "PERSON_EDIT, IMMUNIZATION_VIEW, PERSON_VIEW"
$wsRights.Range("E21").Value = "IMMUNIZATION_VIEW, PERSON_VIEW"
$wsRights.Range("E22").Value = "PERSON_DELETE, IMMUNIZATION_VIEW, VISIT_DELETE, PERSON_VIEW"
$wsRights.Range("E25").Value = "VISIT_DELETE, PERSON_VIEW"
$wsRights.Range("E27").Value = "PERSON_EDIT, PERSON_VIEW"
$wsRights.Range("E31").Value = "ADDITIONAL_TEST_DELETE, PATHOGEN_TEST_DELETE, SAMPLE_VIEW, ADDITIONAL_TEST_VIEW"
$wsRights.Range("E40").Value = "SAMPLE_VIEW, ADDITIONAL_TEST_VIEW"
$wsRights.Range("E41").Value = "SAMPLE_VIEW, ADDITIONAL_TEST_VIEW"
$wsRights.Range("E42").Value = "SAMPLE_VIEW, ADDITIONAL_TEST_VIEW"
$wsRights.Range("E44").Value = "CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E45").Value = "PERSON_EDIT, CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E46").Value = "CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E47").Value = "DOCUMENT_VIEW, VISIT_DELETE, TASK_DELETE, PATHOGEN_TEST_DELETE, PERSON_VIEW, SAMPLE_DELETE, ADDITIONAL_TEST_DELETE, PERSON_DELETE, TASK_VIEW, CONTACT_VIEW, DOCUMENT_DELETE, SAMPLE_VIEW, CASE_VIEW, ADDITIONAL_TEST_VIEW"
$wsRights.Range("E48").Value = "CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E49").Value = "CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E50").Value = "PERSON_EDIT, CONTACT_EDIT, CONTACT_VIEW, PERSON_VIEW, CASE_CREATE, CASE_VIEW"
$wsRights.Range("E51").Value = "PERSON_EDIT, CONTACT_EDIT, CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E52").Value = "PERSON_EDIT, CONTACT_EDIT, CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E53").Value = "PERSON_EDIT, CONTACT_EDIT, CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E63").Value = "TASK_EDIT, TASK_VIEW"
$wsRights.Range("E65").Value = "EVENT_VIEW, DOCUMENT_VIEW, DOCUMENT_DELETE"
$wsRights.Range("E71").Value = "DOCUMENT_VIEW, ACTION_DELETE, VISIT_DELETE, TASK_DELETE, PATHOGEN_TEST_DELETE, PERSON_VIEW, SAMPLE_DELETE, ADDITIONAL_TEST_DELETE, EVENT_VIEW, EVENTPARTICIPANT_DELETE, PERSON_DELETE, TASK_VIEW, DOCUMENT_DELETE, EVENTPARTICIPANT_VIEW, SAMPLE_VIEW, ADDITIONAL_TEST_VIEW"
$wsRights.Range("E74").Value = "EVENT_VIEW, EVENT_EDIT"
$wsRights.Range("E75").Value = "EVENT_VIEW, EVENT_EDIT"
$wsRights.Range("E76").Value = "EVENT_VIEW, PERSON_VIEW"
$wsRights.Range("E77").Value = "EVENT_VIEW, PERSON_VIEW, EVENTPARTICIPANT_VIEW"
$wsRights.Range("E78").Value = "EVENT_VIEW, PERSON_EDIT, PERSON_VIEW, EVENTPARTICIPANT_VIEW"
$wsRights.Range("E79").Value = "EVENT_VIEW, PERSON_VIEW, EVENTPARTICIPANT_VIEW"
$wsRights.Range("E80").Value = "SAMPLE_DELETE, ADDITIONAL_TEST_DELETE, EVENT_VIEW, PERSON_DELETE, VISIT_DELETE, PATHOGEN_TEST_DELETE, PERSON_VIEW, EVENTPARTICIPANT_VIEW, SAMPLE_VIEW, ADDITIONAL_TEST_VIEW"
$wsRights.Range("E81").Value = "EVENT_VIEW, PERSON_VIEW, EVENTPARTICIPANT_VIEW"
$wsRights.Range("E82").Value = "EVENT_VIEW, PERSON_EDIT, EVENTPARTICIPANT_EDIT, PERSON_VIEW, EVENTPARTICIPANT_VIEW"
$wsRights.Range("E87").Value = "EVENT_VIEW, EVENT_EDIT"
$wsRights.Range("E104").Value = "CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E105").Value = "DASHBOARD_CONTACT_VIEW, CONTACT_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E109").Value = "THERAPY_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E110").Value = "THERAPY_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E111").Value = "THERAPY_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E112").Value = "THERAPY_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E113").Value = "THERAPY_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E114").Value = "THERAPY_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E115").Value = "THERAPY_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E116").Value = "THERAPY_VIEW, CLINICAL_COURSE_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E117").Value = "THERAPY_VIEW, CLINICAL_COURSE_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E118").Value = "THERAPY_VIEW, CLINICAL_COURSE_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E119").Value = "THERAPY_VIEW, CLINICAL_COURSE_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E121").Value = "PORT_HEALTH_INFO_VIEW, PERSON_VIEW, CASE_VIEW"
$wsRights.Range("E134").Value = "CAMPAIGN_FORM_DATA_DELETE, CAMPAIGN_FORM_DATA_VIEW, CAMPAIGN_VIEW"
$wsRights.Range("E142").Value = "TRAVEL_ENTRY_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS, PERSON_VIEW"
$wsRights.Range("E143").Value = "PERSON_EDIT, TRAVEL_ENTRY_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS, PERSON_VIEW"
$wsRights.Range("E144").Value = "TRAVEL_ENTRY_VIEW, TRAVEL_ENTRY_MANAGEMENT_ACCESS, PERSON_VIEW"
$wsRights.Range("E145").Value = "DOCUMENT_VIEW, PERSON_DELETE, TASK_VIEW, VISIT_DELETE, TRAVEL_ENTRY_VIEW, TASK_DELETE, DOCUMENT_DELETE, TRAVEL_ENTRY_MANAGEMENT_ACCESS, PERSON_VIEW"
$wsRights.Range("E164").Value = "IMMUNIZATION_DELETE, SAMPLE_EDIT, EVENT_VIEW, PERSON_EDIT, PERSON_DELETE, IMMUNIZATION_VIEW, CONTACT_VIEW, CASE_CREATE, EVENTPARTICIPANT_VIEW, SAMPLE_CREATE, CASE_EDIT, IMMUNIZATION_EDIT, IMMUNIZATION_CREATE, VISIT_DELETE, PATHOGEN_TEST_DELETE, EVENTPARTICIPANT_EDIT, PERSON_VIEW, PATHOGEN_TEST_CREATE, CONTACT_CREATE, CONTACT_EDIT, EVENTPARTICIPANT_CREATE, EVENT_EDIT, EVENT_CREATE, CASE_VIEW, SAMPLE_VIEW, EXTERNAL_MESSAGE_VIEW, PATHOGEN_TEST_EDIT"

# --- Sheet "About": bump SORMAS Version ---
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A2").Value = "1.0.0"

Write-Host "Applied all cell updates."
